$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text in E8 from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Select E8 as the active cell (reflected in sheetView selection)
$ws.Activate()
$ws.Range("E8").Select()
